$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format first so numeric-looking
# strings like "243.03" / "0.3141" are stored verbatim (not coerced to
# floating point numbers, which would corrupt the formatting/precision).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.362.80"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.878.97"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "0.7196"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "243.03"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.07992"
$ws.Range("E8").Value = "  +2.62%  "
$ws.Range("D9").Value = "0.3141"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "24.88"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").Value = "0.08159"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").Value = "1.858.06"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "94.81"
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").Value = "5.226"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "0.7100"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "6.412"
$ws.Range("E16").Value = "  +5.25%  "
$ws.Range("D17").Value = "0.000008464"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "29.366.54"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "247.95"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").Value = "13.29"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "2.125.54"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "7.744"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "0.1603"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "162.85"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").Value = "9.058"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "18.85"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").Value = "1.505"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "4.412"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "4.286"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "1.216"
$ws.Range("E32").Value = "  -5.71%  "
$ws.Range("D33").Value = "0.05342"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").Value = "1.937"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "0.7570"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "1.178"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "0.01881"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").Value = "1.266.35"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("D40").Value = "2.761"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").Value = "6.434"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").Value = "113.30"
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("D43").Value = "0.9075"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("D44").Value = "74.41"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  +4.62%  "
$ws.Range("D47").Value = "2.023.40"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "1.800"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "0.5198"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "9.485"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "0.4342"
$ws.Range("E51").Value = "  +0.01%  "
